$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.458.45"
$ws.Range("E2").Value = "  +0.78%  "
$ws.Range("D3").Value = "2.983.87"
$ws.Range("E3").Value = "  +1.48%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "381.99"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.34"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.546"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.74%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.591"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.69"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.15%  "
$ws.Range("E11").Value = "  -0.82%  "
$ws.Range("E12").Value = "  +1.16%  "
$ws.Range("D13").Value = "3.454.73"
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("E14").Value = "  +2.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.78"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.85%  "
$ws.Range("D16").Value = "2.996.81"
$ws.Range("E16").Value = "  +2.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "11.13"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.996"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").Value = "51.417.81"
$ws.Range("E19").Value = "  +0.79%  "
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("E21").Value = "  +1.54%  "
$ws.Range("E22").Value = "  +0.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.49"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.24"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.21"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.83"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.38"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.94%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").Value = "  +2.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.06"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.89%  "
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.30"
$ws.Range("D32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.62"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "51.57"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.60%  "
$ws.Range("E35").Value = "  +0.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0439"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.58%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.27"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +3.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.77"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.117"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.66%  "
$ws.Range("E41").Value = "  +2.97%  "
$ws.Range("E42").Value = "  +2.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "124.66"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.72"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +11.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.47"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.56%  "
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.39"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.67%  "
$ws.Range("E48").Value = "  -0.85%  "
$ws.Range("D49").Value = "2.027.89"
$ws.Range("E49").Value = "  +1.93%  "
$ws.Range("E50").Value = "  +16.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0329"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.66%  "
